# ============================================================================
# AFDP-6836 FOIA - Fundamental FOIA Changes - Request Created vs. Received Date
# Rewrites the 'Next Possible Queues Rules' decision table on Sheet1:
#   - adds a new 'Appeal' queue step and a 'Default deny queue' action column (H)
#   - replaces the old Suspend/Delete workflow with a deniedFlag-driven Approve/
#     General Counsel/Billing -> 'Approve (Deny)' release path
#   - extends the table by one row (new row 33)
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Prepare layout: bring column H into existence (mirrors column G's
#        formatting on every row of the rule table) and append row 33 by
#        cloning row 32's formatting. ---
foreach ($r in 14,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32) {
    $ws.Range("G$r").Copy()
    $ws.Range("H$r").PasteSpecial(-4122)
}
$ws.Range("A32:H32").Copy()
$ws.Range("A33:H33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Seed the brand-new shared strings in the exact order they are first
#        introduced by the edited table, so the rebuilt table reads naturally. ---
$ws.Range("E26").Value = 'Release,Hold,Fulfill'
$ws.Range("E27").Value = 'General Counsel,Hold,Fulfill'
$ws.Range("E28").Value = 'Billing,Hold,Fulfill'
$ws.Range("E25").Value = 'Release,Fulfill'
$ws.Range("B23").Value = 'Appeal Next Queues'
$ws.Range("C23").Value = 'Appeal'
$ws.Range("E24").Value = 'Hold,Approve'
$ws.Range("H19").Value = '$model.setDefaultDenyQueue($param);'
$ws.Range("H20").Value = 'Default deny queue'
$ws.Range("H21").Value = '"Approve"'
$ws.Range("E22").Value = 'Fulfill,Approve'
$ws.Range("E21").Value = 'Fulfill,Hold,Approve'
$ws.Range("B25").Value = 'Approve to Release (Deny)'
$ws.Range("B32").Value = 'Billing to Approve (Deny)'
$ws.Range("E29").Value = 'Approve,Fulfill'
$ws.Range("E33").Value = 'Release,Fulfill,Approve'
$ws.Range("E31").Value = 'Release,Hold,Fulfill,Approve'
$ws.Range("E30").Value = 'Billing,Hold,Fulfill,Approve'
$ws.Range("B29").Value = 'General Counsel to Approve (Deny)'
$ws.Range("D25").Value = 'deniedFlag'
$ws.Range("D26").Value = '!deniedFlag && !litigationFlag && feeWaiverFlag'
$ws.Range("D27").Value = '!deniedFlag && litigationFlag'
$ws.Range("D28").Value = '!deniedFlag && !litigationFlag && !feeWaiverFlag'
$ws.Range("D30").Value = '!deniedFlag && !feeWaiverFlag'
$ws.Range("D31").Value = '!deniedFlag && feeWaiverFlag'
$ws.Range("D33").Value = '!deniedFlag'

# --- 3. Write the full, final content of the rule table (rows 19-33). ---
$ws.Range("C19").Value = '((CaseFile)$model.getBusinessObject()).getQueue() != null && ((CaseFile)$model.getBusinessObject()).getQueue().getName().equals("$param")'
$ws.Range("D19").Value = 'eval(evalSpring("$param", ((FOIARequest) $model.getBusinessObject())))'
$ws.Range("E19").Value = 'addQueues("$param", $model);'
$ws.Range("F19").Value = '$model.setDefaultNextQueue("$param");'
$ws.Range("G19").Value = '$model.setDefaultReturnQueue($param);'
$ws.Range("H19").Value = '$model.setDefaultDenyQueue($param);'
$ws.Range("A20").Value = 'If different rules apply to the same document, the last rule wins.  Be sure to put default rules first, and specific rules later.'
$ws.Range("B20").Value = 'Rule Name'
$ws.Range("C20").Value = 'Current Queue Name'
$ws.Range("D20").Value = 'When Expression is True'
$ws.Range("E20").Value = 'List of possible next queues'
$ws.Range("F20").Value = 'Default next queue'
$ws.Range("G20").Value = 'Default return queue'
$ws.Range("H20").Value = 'Default deny queue'
$ws.Range("B21").Value = 'Intake Next Queues'
$ws.Range("C21").Value = 'Intake'
$ws.Range("E21").Value = 'Fulfill,Hold,Approve'
$ws.Range("F21").Value = 'Fulfill'
$ws.Range("G21").Value = 'null'
$ws.Range("H21").Value = '"Approve"'
$ws.Range("B22").Value = 'Hold Next Queues'
$ws.Range("C22").Value = 'Hold'
$ws.Range("E22").Value = 'Fulfill,Approve'
$ws.Range("F22").Value = 'Fulfill'
$ws.Range("G22").Value = 'null'
$ws.Range("H22").Value = '"Approve"'
$ws.Range("B23").Value = 'Appeal Next Queues'
$ws.Range("C23").Value = 'Appeal'
$ws.Range("E23").Value = 'Fulfill,Approve'
$ws.Range("F23").Value = 'Fulfill'
$ws.Range("G23").Value = 'null'
$ws.Range("H23").Value = '"Approve"'
$ws.Range("B24").Value = 'Fulfill Next Queues'
$ws.Range("C24").Value = 'Fulfill'
$ws.Range("E24").Value = 'Hold,Approve'
$ws.Range("F24").Value = 'Approve'
$ws.Range("G24").Value = 'null'
$ws.Range("H24").Value = '"Approve"'
$ws.Range("B25").Value = 'Approve to Release (Deny)'
$ws.Range("C25").Value = 'Approve'
$ws.Range("D25").Value = 'deniedFlag'
$ws.Range("E25").Value = 'Release,Fulfill'
$ws.Range("F25").Value = 'Release'
$ws.Range("G25").Value = '"Fulfill"'
$ws.Range("H25").Value = 'null'
$ws.Range("B26").Value = 'Approve – to Release Rules'
$ws.Range("C26").Value = 'Approve'
$ws.Range("D26").Value = '!deniedFlag && !litigationFlag && feeWaiverFlag'
$ws.Range("E26").Value = 'Release,Hold,Fulfill'
$ws.Range("F26").Value = 'Release'
$ws.Range("G26").Value = '"Fulfill"'
$ws.Range("H26").Value = 'null'
$ws.Range("B27").Value = 'Approve to GC'
$ws.Range("C27").Value = 'Approve'
$ws.Range("D27").Value = '!deniedFlag && litigationFlag'
$ws.Range("E27").Value = 'General Counsel,Hold,Fulfill'
$ws.Range("F27").Value = 'General Counsel'
$ws.Range("G27").Value = '"Fulfill"'
$ws.Range("H27").Value = 'null'
$ws.Range("B28").Value = 'Approve to Billing'
$ws.Range("C28").Value = 'Approve'
$ws.Range("D28").Value = '!deniedFlag && !litigationFlag && !feeWaiverFlag'
$ws.Range("E28").Value = 'Billing,Hold,Fulfill'
$ws.Range("F28").Value = 'Billing'
$ws.Range("G28").Value = '"Fulfill"'
$ws.Range("H28").Value = 'null'
$ws.Range("B29").Value = 'General Counsel to Approve (Deny)'
$ws.Range("C29").Value = 'General Counsel'
$ws.Range("D29").Value = 'deniedFlag'
$ws.Range("E29").Value = 'Approve,Fulfill'
$ws.Range("F29").Value = 'Approve'
$ws.Range("G29").Value = '"Fulfill"'
$ws.Range("H29").Value = '"Approve"'
$ws.Range("B30").Value = 'General Counsel to Billing'
$ws.Range("C30").Value = 'General Counsel'
$ws.Range("D30").Value = '!deniedFlag && !feeWaiverFlag'
$ws.Range("E30").Value = 'Billing,Hold,Fulfill,Approve'
$ws.Range("F30").Value = 'Billing'
$ws.Range("G30").Value = '"Fulfill"'
$ws.Range("H30").Value = '"Approve"'
$ws.Range("B31").Value = 'General Counsel to Release'
$ws.Range("C31").Value = 'General Counsel'
$ws.Range("D31").Value = '!deniedFlag && feeWaiverFlag'
$ws.Range("E31").Value = 'Release,Hold,Fulfill,Approve'
$ws.Range("F31").Value = 'Release'
$ws.Range("G31").Value = '"Fulfill"'
$ws.Range("H31").Value = '"Approve"'
$ws.Range("B32").Value = 'Billing to Approve (Deny)'
$ws.Range("C32").Value = 'Billing'
$ws.Range("D32").Value = 'deniedFlag'
$ws.Range("E32").Value = 'Approve,Fulfill'
$ws.Range("F32").Value = 'Approve'
$ws.Range("G32").Value = '"Fulfill"'
$ws.Range("H32").Value = '"Approve"'
$ws.Range("B33").Value = 'Billing to Release'
$ws.Range("C33").Value = 'Billing'
$ws.Range("D33").Value = '!deniedFlag'
$ws.Range("E33").Value = 'Release,Fulfill,Approve'
$ws.Range("F33").Value = 'Release'
$ws.Range("G33").Value = '"Fulfill"'
$ws.Range("H33").Value = '"Approve"'

# --- 4. Propagate the header row's ACTION label and Sequential/true banner
#        into the newly-created column H. ---
$ws.Range("H14").Value = 'true'
$ws.Range("H17").Value = 'ACTION'

# --- 5. Column B/C widened/narrowed to fit the new, longer rule names. ---
$ws.Columns("B").ColumnWidth = 33.7265625
$ws.Columns("C").ColumnWidth = 98.81640625

